$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D/E cells being updated so numeric-looking strings
# (e.g. "591.43", "0.525") are stored as text, matching the source data
# (all Price/Volume cells in this sheet are plain text, not numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = '67.617.15'
$ws.Range("E2").Value = '  +1.04%  '
$ws.Range("D3").Value = '2.529.42'
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '591.43'
$ws.Range("E5").Value = '  +0.53%  '
$ws.Range("D6").Value = '171.13'
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '0.525'
$ws.Range("E8").Value = '  -0.52%  '
$ws.Range("D9").Value = '2.527.77'
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("E10").Value = '  +0.60%  '
$ws.Range("E11").Value = '  +1.20%  '
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("D13").Value = '0.342'
$ws.Range("E13").Value = '  +0.21%  '
$ws.Range("D14").Value = '26.40'
$ws.Range("E14").Value = '  -0.48%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '0.0000177'
$ws.Range("E15").Value = '  +0.94%  '
$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = '2.926.78'
$ws.Range("E16").Value = '  -2.60%  '
$ws.Range("D17").Value = '67.442.14'
$ws.Range("E17").Value = '  +0.91%  '
$ws.Range("D18").Value = '2.516.49'
$ws.Range("E18").Value = '  -0.94%  '
$ws.Range("D19").Value = '11.83'
$ws.Range("E19").Value = '  +4.88%  '
$ws.Range("D20").Value = '7.87'
$ws.Range("E20").Value = '  -1.20%  '
$ws.Range("D21").Value = '368.14'
$ws.Range("E21").Value = '  +4.06%  '
$ws.Range("E22").Value = '  -0.17%  '
$ws.Range("D23").Value = '4.58'
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").Value = '71.66'
$ws.Range("E24").Value = '  +2.93%  '
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("D26").Value = '1.91'
$ws.Range("E26").Value = '  -3.04%  '
$ws.Range("D27").Value = '9.96'
$ws.Range("E27").Value = '  -0.29%  '
$ws.Range("E28").Value = '  -0.34%  '
$ws.Range("D29").Value = '2.640.89'
$ws.Range("E29").Value = '  -1.06%  '
$ws.Range("D30").Value = '0.0₃0961'
$ws.Range("E30").Value = '  -1.50%  '
$ws.Range("D31").Value = '537.18'
$ws.Range("E31").Value = '  +1.00%  '
$ws.Range("E33").Value = '  -0.92%  '
$ws.Range("E34").Value = '  +2.50%  '
$ws.Range("E35").Value = '  -1.16%  '
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("D37").Value = '158.58'
$ws.Range("E37").Value = '  +0.65%  '
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").Value = '19.11'
$ws.Range("E38").Value = '  +3.08%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = '1.42'
$ws.Range("E39").Value = '  -1.87%  '
$ws.Range("E40").Value = '  +1.02%  '
$ws.Range("D41").Value = '0.350'
$ws.Range("E41").Value = '  -1.11%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D42").Value = '5.12'
$ws.Range("E42").Value = '  +0.90%  '
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").Value = '1.78'
$ws.Range("E43").Value = '  -0.15%  '
$ws.Range("E44").Value = '  +0.30%  '
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").Value = '0.0₆0284'
$ws.Range("E46").Value = '  +3.04%  '
$ws.Range("D47").Value = '146.46'
$ws.Range("E47").Value = '  -1.45%  '
$ws.Range("E48").Value = '  -0.55%  '
$ws.Range("D49").Value = '3.70'
$ws.Range("E49").Value = '  +0.91%  '
$ws.Range("D50").Value = '1.72'
$ws.Range("E50").Value = '  +2.09%  '
$ws.Range("D51").Value = '0.0749'
$ws.Range("E51").Value = '  -0.83%  '
